$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 3
    7  = 1
    8  = 2
    9  = 2
    10 = 2
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 3
    20 = 5
    21 = 0
    22 = 4
    23 = 3
    24 = 3
    25 = 5
    26 = 9
    27 = 2
    28 = 7
    29 = 0
    30 = 5
    31 = 2
    32 = 4
    33 = 2
    34 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
